# Zadání.xlsx — "přidání placeholderů do vyhledávače"
#
# Row 29 (vyhledávací formulář úkol) je rozpracovaný -> oznacit jako WIP a
# zvyraznit zlute.
# Row 30 (tlačítko pro výpočet z tabulky) je hotovo -> oznacit jako ano.
# Zavisle souctove/countif vzorce (J3:J7) a F30 se prepocitaji automaticky.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task "Využijte formulář pro práci s daty ... vyhledání a výpis podle
# kritéria" is now work-in-progress: value -> "WIP" + yellow highlight.
$ws.Range("E29").Value = "WIP"
$ws.Range("E29").Interior.Color = 65535

# Task "Formulář bude obsahovat také tlačítko pro libovolný výpočet
# z tabulky" is now done: value -> "ano".
$ws.Range("E30").Value = "ano"

# Match the author's last selection in the sheet.
$ws.Range("G29").Select()
